$wb = $excel.ActiveWorkbook

# --- addListItem sheet: update the escalation user name ---
$wsAddListItem = $wb.Worksheets.Item("addListItem")
$wsAddListItem.Range("A2").Value = "UserEscB"

# --- createUser sheet: bump the test user number ---
$wsCreateUser = $wb.Worksheets.Item("createUser")
$wsCreateUser.Range("A2").Value = 159

# --- Make createUser the active/selected sheet (tab) ---
$wsCreateUser.Activate()
